# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.209.87'
$ws.Range('E2').Value = '  +7.25%  '
$ws.Range('D3').Value = '2.678.88'
$ws.Range('E3').Value = '  +10.06%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '514.24'
$ws.Range('E5').Value = '  +4.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.93'
$ws.Range('E6').Value = '  +2.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.620'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').Value = '2.673.65'
$ws.Range('E9').Value = '  +9.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.15'
$ws.Range('E10').Value = '  +9.21%  '
$ws.Range('E11').Value = '  +5.83%  '
$ws.Range('E12').Value = '  +4.35%  '
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').Value = '3.131.01'
$ws.Range('E14').Value = '  +9.61%  '
$ws.Range('D15').Value = '61.266.13'
$ws.Range('E15').Value = '  +7.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.41'
$ws.Range('E16').Value = '  +7.74%  '
$ws.Range('E17').Value = '  +5.05%  '
$ws.Range('D18').Value = '2.669.70'
$ws.Range('E18').Value = '  +9.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.87'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '353.87'
$ws.Range('E20').Value = '  +7.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.60'
$ws.Range('E21').Value = '  +6.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.21'
$ws.Range('E22').Value = '  +4.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.89'
$ws.Range('E24').Value = '  +4.85%  '
$ws.Range('E25').Value = '  +3.80%  '
$ws.Range('D26').Value = '2.780.72'
$ws.Range('E26').Value = '  +9.12%  '
$ws.Range('E27').Value = '  +4.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +10.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.61'
$ws.Range('E30').Value = '  +2.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.74'
$ws.Range('E32').Value = '  +4.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '157.02'
$ws.Range('E33').Value = '  +4.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.59'
$ws.Range('E34').Value = '  +4.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.77'
$ws.Range('E35').Value = '  +8.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.14'
$ws.Range('E36').Value = '  +11.50%  '
$ws.Range('E37').Value = '  +7.44%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.55'
$ws.Range('E38').Value = '  +12.01%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.886'
$ws.Range('E39').Value = '  +3.34%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '307.87'
$ws.Range('E40').Value = '  +15.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.82'
$ws.Range('E41').Value = '  +8.20%  '
$ws.Range('B42').Value = 'SuiNetwork'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.847'
$ws.Range('E42').Value = '  +32.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '35.80'
$ws.Range('E43').Value = '  +4.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.648'
$ws.Range('E44').Value = '  +8.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0584'
$ws.Range('E45').Value = '  +8.54%  '
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.32'
$ws.Range('E47').Value = '  +15.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.998'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.02'
$ws.Range('E49').Value = '  +7.10%  '
$ws.Range('E50').Value = '  +4.01%  '
$ws.Range('D51').Value = '2.034.73'
$ws.Range('E51').Value = '  +9.21%  '
